$wb = $excel.ActiveWorkbook

# "customcolumns" sheet (sheet1.xml) gains new data in rows 12-15 (columns A & B).
$ws = $wb.Worksheets.Item("customcolumns")

# Fill in the new rows in the same order they were authored so that the
# shared-string table is built up with matching indices/order.
$ws.Range("A12").Value = "sub_category_src_code"
$ws.Range("A13").Value = "sub_category_src_key"
$ws.Range("B12").Value = "dwh_f_incident"
$ws.Range("B13").Value = "dwh_f_incident"
$ws.Range("B14").Value = "dwh_d_incident"
$ws.Range("B15").Value = "dwh_d_problem"
$ws.Range("A14").Value = "short_description"
$ws.Range("A15").Value = "short_description"

# Move the active sheet/selection to "customcolumns" (was "PLP_Updated").
# Activating this sheet automatically clears the "tabSelected" flag on
# "PLP_Updated" (sheet3.xml), whose own A3 selection is left untouched.
$ws.Activate()
$ws.Range("A14").Select()
